$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.591.57'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.65%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.922.77'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.88%  '

$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.67%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.14%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4710'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.49%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2897'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.49%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06789'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.26%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '105.27'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.97%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '18.41'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.74%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.904.52'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.93%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07700'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.21%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.293'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6743'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +8.08%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '290.50'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.13%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.611.11'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.75%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000007615'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +3.43%  '

$ws.Range("E19").Value = '  -0.07%  '

$ws.Range("E20").Value = '  +2.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.160.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.457'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +9.48%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.001'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.340'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.69%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.410'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.84%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.10'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.86%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.97'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +9.15%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.124'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +9.83%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.1080'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.71%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.378'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.97%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.185'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.74%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.146'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +8.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05062'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +3.27%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7422'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +2.30%  '

$ws.Range("E35").Value = '  +3.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02088'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.750'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.693'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +1.43%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.061'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +4.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '111.20'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.51%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8817'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.72%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4352'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +7.85%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.891'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.09%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.000'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.13%  '

$ws.Range("E45").Value = '  +3.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.243'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +3.06%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '9.304'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.86%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '47.99'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +18.60%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1234'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.03%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '34.97'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.22%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4041'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +8.86%  '
